$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.016657
$ws.Range("N2").Value = 0.033314
$ws.Range("Q2").Value = 0.007100151744333334
$ws.Range("R2").Value = 0.04260091046600001
